$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "2025/12/03 05:00"
$ws.Range("B14").Value = "43,759位本"
$ws.Range("C14").Value = "97位 広告・宣伝 (本)"
$ws.Range("D14").Value = "175位商業デザイン"
$ws.Range("E14").Value = "2,145位ビジネス実用本"
$ws.Range("F14").Value = "-"
$ws.Range("G14").Value = "-"
